$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Turn the two "Hsd" formula cells (I2:I3) into plain literal date values
#    *before* we delete column H - otherwise the formula (which refers to H)
#    would break into #REF! once H disappears.
# ---------------------------------------------------------------------------
$ws.Range("I2").NumberFormat = "mm-dd-yy"
$ws.Range("I2").Value = 44184
$ws.Range("I3").NumberFormat = "mm-dd-yy"
$ws.Range("I3").Value = 44184

# ---------------------------------------------------------------------------
# 2) Drop the dropdown validations that lived on "Don vi" (F) and
#    "Nha phan phoi" (J) - the target sheet no longer restricts those cells.
# ---------------------------------------------------------------------------
$ws.Range("F2:F3").Validation.Delete()
$ws.Range("J2:J3").Validation.Delete()

# ---------------------------------------------------------------------------
# 3) Clear out the sample values that used to live in those two columns.
# ---------------------------------------------------------------------------
$ws.Range("F2:F3").ClearContents()
$ws.Range("J2:J3").ClearContents()

# ---------------------------------------------------------------------------
# 4) Remove the whole "Hsd(thang)" column (old H) - everything to its right
#    (old I..T) shifts one column to the left (new H..S).
# ---------------------------------------------------------------------------
$ws.Columns("H").Delete()

# ---------------------------------------------------------------------------
# 5) Re-apply number formats, now using the *post-shift* column letters.
# ---------------------------------------------------------------------------

# Header row - give every header the same number format as the data below it.
$ws.Range("C1").NumberFormat = "#,##0_);[Red](#,##0)"
$ws.Range("D1").NumberFormat = "0_);[Red](0)"
$ws.Range("E1").NumberFormat = "#,##0_);[Red](#,##0)"
$ws.Range("G1").NumberFormat = "#,##0_);[Red](#,##0)"
$ws.Range("H1").NumberFormat = "mm-dd-yy"
$ws.Range("K1").NumberFormat = "0.0_ "
$ws.Range("L1").NumberFormat = "@"
$ws.Range("N1").NumberFormat = "@"
$ws.Range("P1").NumberFormat = "@"
$ws.Range("R1").NumberFormat = "@"

# Data rows.
$ws.Range("C2:C3").NumberFormat = "#,##0_);[Red](#,##0)"
$ws.Range("D2:D3").NumberFormat = "0_);[Red](0)"
$ws.Range("E2:E3").NumberFormat = "#,##0_);[Red](#,##0)"
$ws.Range("G2:G3").NumberFormat = "#,##0_);[Red](#,##0)"
$ws.Range("H2:H3").NumberFormat = "mm-dd-yy"
$ws.Range("K2:K3").NumberFormat = "0.0_ "

# ---------------------------------------------------------------------------
# 6) "Gia ban" (G2:G3) drops its left border now that "Don vi" (F) is blank.
# ---------------------------------------------------------------------------
$ws.Range("G2:G3").Borders.Item(7).LineStyle = -4142

Write-Host "edit complete"
